$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.667.45'
$ws.Range('E2').Value = '  +1.43%  '

$ws.Range('D3').Value = '3.628.96'
$ws.Range('E3').Value = '  +3.53%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').Value = '''605.40'
$ws.Range('E5').Value = '  -0.06%  '

$ws.Range('D6').Value = '''198.59'
$ws.Range('E6').Value = '  +1.95%  '

$ws.Range('E7').Value = '  +0.22%  '

$ws.Range('E8').Value = '  +0.10%  '

$ws.Range('E9').Value = '  +9.63%  '

$ws.Range('D10').Value = '''0.646'
$ws.Range('E10').Value = '  -0.62%  '

$ws.Range('D11').Value = '''53.82'
$ws.Range('E11').Value = '  +0.58%  '

$ws.Range('D12').Value = '''0.0000305'
$ws.Range('E12').Value = '  +1.43%  '

$ws.Range('D13').Value = '''9.54'
$ws.Range('E13').Value = '  +0.23%  '

$ws.Range('D14').Value = '4.206.59'
$ws.Range('E14').Value = '  +3.51%  '

$ws.Range('D15').Value = '''671.37'
$ws.Range('E15').Value = '  +13.06%  '

$ws.Range('D16').Value = '''13.02'
$ws.Range('E16').Value = '  +2.10%  '

$ws.Range('D17').Value = '70.862.52'
$ws.Range('E17').Value = '  +1.47%  '

$ws.Range('D18').Value = '3.636.69'
$ws.Range('E18').Value = '  +3.73%  '

$ws.Range('D19').Value = '''19.06'
$ws.Range('E19').Value = '  -0.36%  '

$ws.Range('E20').Value = '  +0.34%  '

$ws.Range('D21').Value = '''0.998'
$ws.Range('E21').Value = '  +0.85%  '

$ws.Range('D22').Value = '''18.64'
$ws.Range('E22').Value = '  +2.21%  '

$ws.Range('D23').Value = '''5.38'
$ws.Range('E23').Value = '  +1.60%  '

$ws.Range('D24').Value = '''105.46'
$ws.Range('E24').Value = '  +3.74%  '

$ws.Range('E25').Value = '  -0.61%  '

$ws.Range('E26').Value = '  -5.31%  '

$ws.Range('D27').Value = '''10.44'
$ws.Range('E27').Value = '  -3.83%  '

$ws.Range('D28').Value = '''9.81'
$ws.Range('E28').Value = '  +2.92%  '

$ws.Range('D29').Value = '''34.04'
$ws.Range('E29').Value = '  +2.46%  '

$ws.Range('D30').Value = '''4.68'
$ws.Range('E30').Value = '  +8.98%  '

$ws.Range('D31').Value = '''7.16'
$ws.Range('E31').Value = '  +1.08%  '

$ws.Range('D32').Value = '''12.19'
$ws.Range('E32').Value = '  -1.73%  '

$ws.Range('E33').Value = '  +0.10%  '

$ws.Range('D34').Value = '''63.41'
$ws.Range('E34').Value = '  +0.45%  '

$ws.Range('D35').Value = '3.969.39'
$ws.Range('E35').Value = '  +6.69%  '

$ws.Range('D36').Value = '0.0₃0864'
$ws.Range('E36').Value = '  +4.64%  '

$ws.Range('E37').Value = '  -0.02%  '

$ws.Range('E38').Value = '  -1.70%  '

$ws.Range('D39').Value = '''505.97'
$ws.Range('E39').Value = '  +4.60%  '

$ws.Range('D40').Value = '''36.76'
$ws.Range('E40').Value = '  +1.03%  '

$ws.Range('D41').Value = '''0.388'
$ws.Range('E41').Value = '  -0.94%  '

$ws.Range('E42').Value = '  -3.33%  '

$ws.Range('E43').Value = '  +2.37%  '

$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '''0.0460'
$ws.Range('E44').Value = '  +1.47%  '

$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').Value = '''3.06'
$ws.Range('E45').Value = '  +8.53%  '

$ws.Range('D46').Value = '''3.50'
$ws.Range('E46').Value = '  +6.70%  '

$ws.Range('E47').Value = '  +0.52%  '

$ws.Range('D48').Value = '''8.67'
$ws.Range('E48').Value = '  +3.24%  '

$ws.Range('E49').Value = '  -0.27%  '

$ws.Range('E50').Value = '  +0.66%  '

$ws.Range('E51').Value = '  +5.08%  '
